$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: fix double space "que  influenciaram" -> "que influenciaram"
# ------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("que  influenciaram", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "que influenciaram", 2)
Write-Host "Change1 found=$found1"

# ------------------------------------------------------------------
# Change 2: insert " e alta demanda de consumidores por aquele tipo de
# alimento" between "chinesa" and ". Ele reparou", then split that
# inserted span into its own run (matching the target OOXML which has
# three separate <w:r> elements with identical formatting).
# ------------------------------------------------------------------
$old2 = "chinesa. Ele reparou"
$new2 = "chinesa e alta demanda de consumidores por aquele tipo de alimento. Ele reparou"
$r2 = $d.Content
$found2 = $r2.Find.Execute($old2, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $new2, 2)
Write-Host "Change2 insert found=$found2"

$r2b = $d.Content
$found2b = $r2b.Find.Execute(" e alta demanda de consumidores por aquele tipo de alimento", `
                              $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Change2 locate found=$found2b Start=$($r2b.Start) End=$($r2b.End)"

$toggle2 = $d.Range($r2b.Start, $r2b.End)
$toggle2.Bold = 1
$toggle2.Bold = 0

# ------------------------------------------------------------------
# Change 3: rewrite the final paragraph.
#   - merge the many small runs that make up the paragraph into one
#   - insert the word " empresarial" before " desta dimensão."
#   - move the "_GoBack" bookmark so that it now sits right after the
#     newly inserted " empresarial" text
# ------------------------------------------------------------------

# Remove the bookmark first so that a single Find/Replace can freely
# span across its old position without Word silently dropping it.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$fullOld3 = "O capital inicial para a criação da empresa veio da venda de dois consultórios que Robison Shiba tinha e através do seu pai, que vendeu um apartamento por um terço do valor para ajudar o filho nesta empreitada. Embora todo o negocio inicialmente tenha riscos, este fato ajudou Shiba a não contrair uma divida no inicio da empreitada, entretanto, não é informado o valor do capital inicial que foi utilizado para uma iniciativa desta dimensão."
$fullNew3 = "O capital inicial para a criação da empresa veio da venda de dois consultórios que Robison Shiba tinha e através do seu pai, que vendeu um apartamento por um terço do valor para ajudar o filho nesta empreitada. Embora todo o negocio inicialmente tenha riscos, este fato ajudou Shiba a não contrair uma divida no inicio da empreitada, entretanto, não é informado o valor do capital inicial que foi utilizado para uma iniciativa empresarial desta dimensão."

$r3 = $d.Content
$found3 = $r3.Find.Execute($fullOld3, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $fullNew3, 2)
Write-Host "Change3 merge found=$found3"

# Locate the newly-added " empresarial" word so we can split it into
# its own run and park the bookmark right after it.
$r3b = $d.Content
$found3b = $r3b.Find.Execute(" empresarial", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
Write-Host "Change3 locate found=$found3b Start=$($r3b.Start) End=$($r3b.End)"

$splitStart = $r3b.Start
$splitEnd = $r3b.End

$toggle3 = $d.Range($splitStart, $splitEnd)
$toggle3.Bold = 1
$toggle3.Bold = 0

$bmRange = $d.Range($splitEnd, $splitEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "DONE"
